$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-09-04 12:19:57"
$overview.Range("G5").Value = "2016-09-04 12:19:57"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$zhcn.Range("H4").Value = "2016-09-04 12:19:53"
$zhcn.Range("H5").Value = "2016-09-04 12:19:53"
$zhcn.Range("K4").Value = "2016-09-04 12:20:19"
$zhcn.Range("K5").Value = "2016-09-04 12:20:19"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K4").Value = "2016-09-04 12:20:26"
$dede.Range("K5").Value = "2016-09-04 12:20:26"
